# Refresh cached Universalis market-price columns (H:N) in the leve-profit sheets.
# Mirrors the scheduled runner's nightly re-pull of currentAveragePrice* / LevePrice* /
# LeveProfit* figures for a handful of leves across the ALC, ARM, CRP, CUL, GSM, LTW and
# WVR crafter tabs. Only the numeric market-data cells are touched -- leve name/item/level
# /EXP/gil/amount/item-id columns (A:G) are left untouched.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 43: Growing Is Knowing
$ws.Range("H43").Value = 1356.7059
$ws.Range("J43").Value = 1356.7059
$ws.Range("L43").Value = 1356.7059
$ws.Range("N43").Value = -1494.7059

# Row 101: Edge of the Arcane
$ws.Range("H101").Value = 1533
$ws.Range("I101").Value = 399.66666
$ws.Range("K101").Value = 1198.99998
$ws.Range("M101").Value = 423.0000199999999

# Row 116: Growing Up
$ws.Range("H116").Value = 16410
$ws.Range("I116").Value = 100000
$ws.Range("J116").Value = 7122.222
$ws.Range("K116").Value = 100000
$ws.Range("L116").Value = 7122.222
$ws.Range("M116").Value = -96558
$ws.Range("N116").Value = -14006.222

# Row 124: Luncheon Bound
$ws.Range("H124").Value = 47247.25
$ws.Range("J124").Value = 47247.25
$ws.Range("L124").Value = 47247.25
$ws.Range("N124").Value = -57067.25

# Row 138: All-night Crafting
$ws.Range("H138").Value = 2856.4666
$ws.Range("I138").Value = 2653.68
$ws.Range("J138").Value = 3870.4
$ws.Range("K138").Value = 7961.039999999999
$ws.Range("L138").Value = 11611.2
$ws.Range("M138").Value = -2821.039999999999
$ws.Range("N138").Value = -21891.2

# Row 141: Remedy for Reason
$ws.Range("H141").Value = 3829.8125
$ws.Range("J141").Value = 4117.9
$ws.Range("L141").Value = 12353.7
$ws.Range("N141").Value = -22713.7

$ws = $wb.Worksheets.Item("ARM")
# Row 32: Ingot We Trust
$ws.Range("H32").Value = 3903.1667
$ws.Range("I32").Value = 2098.1667
$ws.Range("K32").Value = 2098.1667
$ws.Range("M32").Value = -1811.1667

# Row 61: Dealing with the Tough Stuff
$ws.Range("H61").Value = 5538
$ws.Range("I61").Value = 3082.8333
$ws.Range("K61").Value = 3082.8333
$ws.Range("M61").Value = -2870.8333

# Row 74: As the Bolt Flies
$ws.Range("H74").Value = 1120.6666
$ws.Range("I74").Value = 1046.8462
$ws.Range("K74").Value = 1046.8462
$ws.Range("M74").Value = -172.8462

# Row 77: Heavy Metal Banned (L)
$ws.Range("H77").Value = 1120.6666
$ws.Range("I77").Value = 1046.8462
$ws.Range("K77").Value = 5234.231
$ws.Range("M77").Value = -866.2309999999998

# Row 102: Smells of Rich Tama-hagane
$ws.Range("H102").Value = 1797.625
$ws.Range("J102").Value = 3287.25
$ws.Range("L102").Value = 3287.25
$ws.Range("N102").Value = -6531.25

# Row 136: Metal with Mettle
$ws.Range("H136").Value = 5538
$ws.Range("I136").Value = 3082.8333
$ws.Range("K136").Value = 9248.499899999999
$ws.Range("M136").Value = -6698.499899999999

$ws = $wb.Worksheets.Item("CRP")
# Row 58: You Do the Heavy Lifting
$ws.Range("H58").Value = 1318806.5
$ws.Range("I58").Value = 2416507.5
$ws.Range("K58").Value = 2416507.5
$ws.Range("M58").Value = -2416304.5

# Row 134: Wood You Be Quiet
$ws.Range("H134").Value = 1799.2106
$ws.Range("I134").Value = 1343.6666
$ws.Range("K134").Value = 4030.9998
$ws.Range("M134").Value = -1495.9998

# Row 136: Turali Quality
$ws.Range("H136").Value = 1318806.5
$ws.Range("I136").Value = 2416507.5
$ws.Range("K136").Value = 7249522.5
$ws.Range("M136").Value = -7246972.5

$ws = $wb.Worksheets.Item("CUL")
# Row 3: Trout Fishing in Limsa
$ws.Range("H3").Value = 3000
$ws.Range("I3").Value = 3000
$ws.Range("K3").Value = 9000
$ws.Range("M3").Value = -8888

# Row 75: Breakfast of Champions
$ws.Range("H75").Value = 949
$ws.Range("J75").Value = 949
$ws.Range("L75").Value = 2847
$ws.Range("N75").Value = -4843

# Row 78: Emerald Soup for the Soul (L)
$ws.Range("H78").Value = 949
$ws.Range("J78").Value = 949
$ws.Range("L78").Value = 8541
$ws.Range("N78").Value = -18525

# Row 106: Herky Jerky
$ws.Range("H106").Value = 0
$ws.Range("J106").Value = 0
$ws.Range("L106").Value = 0
$ws.Range("N106").Value = $null

# Row 129: Comfort Food
$ws.Range("H129").Value = 61351.332
$ws.Range("J129").Value = 104685.71
$ws.Range("L129").Value = 314057.13
$ws.Range("N129").Value = -324057.13

# Row 131: The Mountain Steeped
$ws.Range("H131").Value = 780.48
$ws.Range("J131").Value = 799.9681
$ws.Range("L131").Value = 2399.9043
$ws.Range("N131").Value = -12479.9043

$ws = $wb.Worksheets.Item("GSM")
# Row 80: Needs More Prayerbell
$ws.Range("H80").Value = 2335.75
$ws.Range("I80").Value = 2281
$ws.Range("J80").Value = 2500
$ws.Range("K80").Value = 2281
$ws.Range("L80").Value = 2500
$ws.Range("M80").Value = -1283
$ws.Range("N80").Value = -4496

# Row 83: With a Noise That Reaches Heaven (L)
$ws.Range("H83").Value = 2335.75
$ws.Range("I83").Value = 2281
$ws.Range("J83").Value = 2500
$ws.Range("K83").Value = 11405
$ws.Range("L83").Value = 12500
$ws.Range("M83").Value = -6413
$ws.Range("N83").Value = -22484

# Row 102: Put the Metal to the Peddle
$ws.Range("H102").Value = 2155.4644
$ws.Range("I102").Value = 2002.1052
$ws.Range("K102").Value = 2002.1052
$ws.Range("M102").Value = -380.1052

$ws = $wb.Worksheets.Item("LTW")
# Row 7: Tan Before the Ban
$ws.Range("H7").Value = 8771.571
$ws.Range("I7").Value = 5902
$ws.Range("J7").Value = 9919.4
$ws.Range("K7").Value = 5902
$ws.Range("L7").Value = 9919.4
$ws.Range("M7").Value = -5790
$ws.Range("N7").Value = -10143.4

# Row 22: Skin off Their Backs
$ws.Range("H22").Value = 2293.111
$ws.Range("I22").Value = 2639.8
$ws.Range("J22").Value = 1859.75
$ws.Range("K22").Value = 2639.8
$ws.Range("L22").Value = 1859.75
$ws.Range("M22").Value = -2344.8
$ws.Range("N22").Value = -2449.75

# Row 27: Fire and Hide
$ws.Range("H27").Value = 2293.111
$ws.Range("I27").Value = 2639.8
$ws.Range("J27").Value = 1859.75
$ws.Range("K27").Value = 2639.8
$ws.Range("L27").Value = 1859.75
$ws.Range("M27").Value = -2532.8
$ws.Range("N27").Value = -2073.75

# Row 55: It's Not a Job, It's a Calling
$ws.Range("H55").Value = 405.66666
$ws.Range("I55").Value = 349.6
$ws.Range("K55").Value = 349.6
$ws.Range("M55").Value = -176.6

# Row 126: Battered Books
$ws.Range("H126").Value = 8771.571
$ws.Range("I126").Value = 5902
$ws.Range("J126").Value = 9919.4
$ws.Range("K126").Value = 17706
$ws.Range("L126").Value = 29758.2
$ws.Range("M126").Value = -15236
$ws.Range("N126").Value = -34698.2

# Row 136: Respect for Br'aax
$ws.Range("H136").Value = 4205.2856
$ws.Range("J136").Value = 5961.875
$ws.Range("L136").Value = 17885.625
$ws.Range("N136").Value = -22985.625

$ws = $wb.Worksheets.Item("WVR")
# Row 113: A Tender Table
$ws.Range("H113").Value = 687.6818
$ws.Range("I113").Value = 431.11765
$ws.Range("J113").Value = 1560
$ws.Range("K113").Value = 1293.35295
$ws.Range("L113").Value = 4680
$ws.Range("M113").Value = 876.64705
$ws.Range("N113").Value = -9020

# Row 122: Heavy Armoire
$ws.Range("H122").Value = 100244.69
$ws.Range("I122").Value = 145534.53
$ws.Range("K122").Value = 436603.59
$ws.Range("M122").Value = -434153.59

# Row 126: A Polished Purchase
$ws.Range("H126").Value = 6877.1333
$ws.Range("I126").Value = 6236.1113
$ws.Range("K126").Value = 18708.3339
$ws.Range("M126").Value = -16238.3339
